# "End of section 4" - add two more NuGet packages to the Nugget section
# of the Cmd sheet, and make the Cmd sheet the active tab (previously the
# Angular sheet was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cmd")

# Insert two new rows right after the existing Nugget entries (row 23)
# and before the blank spacer row that preceded the VERSIONS section.
$ws.Rows("24:25").Insert()

$ws.Range("B24").Value = "System.IdentityModel.Tokens.Jwt"
$ws.Range("B25").Value = "Microsoft.AspNetCore.Authentication.JwtBearer"

# Column B's best-fit width shrinks now that the font metrics were
# recalculated for the sheet.
$ws.Columns("B").ColumnWidth = 46.67

# Update the selection to the new last row.
$ws.Range("C32").Select() | Out-Null

# Make "Cmd" the active/selected sheet (moves tabSelected from Angular).
$ws.Activate() | Out-Null

Write-Output "done"
